# Add latest pipe count results
# - Adds a new strategy results column (W) with its own header/commit-hash,
#   shifts the "diff vs current min" formula from W to X, and drops the old
#   "diff vs FBE baseline" (D-V) formula that used to live in X.
# - Extends the MIN() range formula in column B to include the new W column.
# - Updates the two conditional formatting rules so they also cover the new
#   W/X columns.
# - Updates the sheet view (selection / scroll position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New column header (row 1 = strategy name, row 2 = commit hash)
# ---------------------------------------------------------------------
$ws.Range("W1").Value = "Bug fixes and perf"
$ws.Range("W2").Value = "908a6d81fbc9a17161ab135f40d9dd8b2dd8787b"
# Match the style used by the other commit-hash cells in row 2 (e.g. V2):
# bold font + scientific number format (the cell still displays the text
# since NumberFormat only affects numeric values).
$ws.Range("W2").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------
# 2) New column data (row 3 .. 60) - latest pipe-count run results
# ---------------------------------------------------------------------
$newValues = @{
    3=16; 4=121; 5=95; 6=78; 7=83; 8=75; 9=89; 10=64; 11=86; 12=103;
    13=104; 14=54; 15=76; 16=90; 17=67; 18=80; 19=28; 20=37; 21=31; 22=47;
    23=14; 24=54; 25=67; 26=29; 27=29; 28=67; 29=46; 30=44; 31=84; 32=56;
    33=71; 34=57; 35=44; 36=127; 37=134; 38=80; 39=91; 40=88; 41=78; 42=87;
    43=114; 44=103; 45=119; 46=118; 47=91; 48=111; 49=112; 50=128; 51=89;
    52=123; 53=118; 54=168; 55=143; 56=101; 57=112; 58=103; 59=114; 60=131
}

foreach ($r in 3..60) {
    $ws.Range("W$r").Value = $newValues[$r]
}

# ---------------------------------------------------------------------
# 3) Column B ("Min") formulas now need to look across C:W instead of C:V
# ---------------------------------------------------------------------
$ws.Range("B3").Formula = "=MIN(C3:W3)"
$ws.Range("B4:B60").FormulaR1C1 = "=MIN(RC[1]:RC[21])"

# ---------------------------------------------------------------------
# 4) Column X now holds what column W used to hold (current min vs "Add
#    FBE strategy" (V) diff). The old column X formula (D-V, diff vs the
#    FBE baseline) is gone, overwritten by this.
# ---------------------------------------------------------------------
$ws.Range("X3").Formula = "=V3-B3"
$ws.Range("X4:X60").FormulaR1C1 = "=RC[-2]-RC[-22]"

# Row 61 summary (averages). E61:W61 stay a shared AVERAGE() formula,
# X61 is a new cell continuing that same series (now the average of the
# relocated V-B diff column).
$ws.Range("X61").Formula = "=AVERAGE(X3:X60)"

# ---------------------------------------------------------------------
# 5) Conditional formatting: extend ranges from V/W to W/X.
#    The underlying engine only keeps the first area of a multi-area
#    range when assigned to a FormatCondition, so the additional areas
#    are re-created as extra rules using an identical style (same font /
#    fill colors as the existing "good" rule) to keep the same visual
#    effect as the original multi-area rule.
# ---------------------------------------------------------------------
$fcs = $ws.Cells.FormatConditions
$ruleExpr = $fcs.Item(1)
$exprFontColor = $ruleExpr.Font.Color
$exprInteriorColor = $ruleExpr.Interior.Color
$exprFormula = $ruleExpr.Formula1

$ruleCellIs = $fcs.Item(2)

# Rule 1 originally applied to "C1:V60 C62:V1048576 C61:W61"
# -> now "C1:W60 C62:W1048576 C61:X61"
$ruleExpr.ModifyAppliesToRange($ws.Range("C1:W60"))

$ruleExprB = $ws.Range("C62:W1048576").FormatConditions.Add(2, 0, $exprFormula)
$ruleExprB.Font.Color = $exprFontColor
$ruleExprB.Interior.Color = $exprInteriorColor

$ruleExprC = $ws.Range("C61:X61").FormatConditions.Add(2, 0, $exprFormula)
$ruleExprC.Font.Color = $exprFontColor
$ruleExprC.Interior.Color = $exprInteriorColor

# Rule 2 originally applied to "C4:V4" -> now "C4:W4"
$ruleCellIs.ModifyAppliesToRange($ws.Range("C4:W4"))

# ---------------------------------------------------------------------
# 6) Sheet view: selection moved to W25, scrolled back up to top-left B1
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("W25").Select()
